# ExtractJobDetails - first commit
# Populates the second job-listing row (row 2) with the full set of
# scraped job-detail columns (Job Ref .. Years) for the "QC inspector"
# listing at E-STAR PRECISION (S) PTE LTD, pushing the previously
# second row (2.pdf / AutoCAD,... ) down to row 3 - matching the
# authoritative ordering produced by the scraping run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: job detail columns F..O for the newly extracted job (132.pdf)
$ws.Range("F2").Value = "MCF-2022-0116373"
$ws.Range("G2").Value = "E-STAR PRECISION (S) PTE LTD"
$ws.Range("H2").Value = "QC inspector "
$ws.Range("I2").Value = "0% skills matched"
$ws.Range("J2").Value = "NA"
$ws.Range("K2").Value = "Permanent"
$ws.Range("L2").Value = "Non-executive"
$ws.Range("M2").Value = "Engineering"
$ws.Range("N2").Value = "2035 BUKIT BATOK STREET 23 659540"
$ws.Range("O2").Value = "1 year exp"

# Standardise the column widths across the data range (A:L narrow data
# columns, M the wider "Industry"/address column) to match the default
# layout written by the newer Excel build that produced this workbook.
$ws.Columns("A:L").ColumnWidth = 8.25
$ws.Columns("M:M").ColumnWidth = 10.42

# Move the active selection to O1, as left by the authoring session.
$ws.Range("O1").Select()
